$d = $word.ActiveDocument

# Insert three new paragraphs right after the first paragraph ("example" / Title)
$anchor = $d.Paragraphs(1).Range
$anchor.InsertParagraphAfter()

$p2 = $d.Paragraphs(2).Range
$p2.Text = "Demo example used for the testing of the ipxact2systemverilog tool."
$d.Paragraphs(2).Style = "FirstParagraph"

$p2End = $d.Paragraphs(2).Range
$p2End.InsertParagraphAfter()

$p3 = $d.Paragraphs(3).Range
$p3.Text = "Base Address"
$d.Paragraphs(3).Style = "DefinitionTerm"

$p3End = $d.Paragraphs(3).Range
$p3End.InsertParagraphAfter()

$p4 = $d.Paragraphs(4).Range
$p4.Text = "0x0"
$d.Paragraphs(4).Style = "Definition"
